$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,16

$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 0.4593233333333334
$data[0,3] = 1.37797
$data[0,4] = 0.015538272766109
$data[0,5] = 0.015538272766109
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 1.664391
$data[0,9] = 4.993173000000001
$data[0,10] = 0.3990511495040125
$data[0,11] = 0.3990511495040125
$data[0,12] = 0.7644936220900002
$data[0,13] = 6.880442598810001
$data[0,14] = 0.00620056560862269
$data[0,15] = 0.006200565608622691

$data[1,0] = 3
$data[1,1] = 1
$data[1,2] = 0.4593233333333334
$data[1,3] = 1.37797
$data[1,4] = 0.015538272766109
$data[1,5] = 0.015538272766109
$data[1,6] = 3
$data[1,7] = 1
$data[1,8] = 0.437958
$data[1,9] = 1.313874
$data[1,10] = 0.1050039584054939
$data[1,11] = 0.1050039584054938
$data[1,12] = 0.20116432842
$data[1,13] = 1.81047895578
$data[1,14] = 0.001631580147225728
$data[1,15] = 0.001631580147225728

$data[2,0] = 3
$data[2,1] = 1
$data[2,2] = 0.4593233333333334
$data[2,3] = 1.37797
$data[2,4] = 0.015538272766109
$data[2,5] = 0.015538272766109
$data[2,6] = 3
$data[2,7] = 1
$data[2,8] = 1.439215333333333
$data[2,9] = 4.317646
$data[2,10] = 0.3450634695516061
$data[2,11] = 0.3450634695516061
$data[2,12] = 0.6610651842911112
$data[2,13] = 5.94958665862
$data[2,14] = 0.005361690311512804
$data[2,15] = 0.005361690311512804

$data[3,0] = 3
$data[3,1] = 1
$data[3,2] = 0.4593233333333334
$data[3,3] = 1.37797
$data[3,4] = 0.015538272766109
$data[3,5] = 0.015538272766109
$data[3,6] = 3
$data[3,7] = 1
$data[3,8] = 0.629307
$data[3,9] = 1.887921
$data[3,10] = 0.1508814225388875
$data[3,11] = 0.1508814225388875
$data[3,12] = 0.28905538893
$data[3,13] = 2.60149850037
$data[3,14] = 0.002344436698747781
$data[3,15] = 0.002344436698747782

$data[4,0] = 3
$data[4,1] = 1
$data[4,2] = 11.94315233333334
$data[4,3] = 35.829457
$data[4,4] = 0.4040203167903319
$data[4,5] = 0.4040203167903319
$data[4,6] = 3
$data[4,7] = 1
$data[4,8] = 1.664391
$data[4,9] = 4.993173000000001
$data[4,10] = 0.3990511495040125
$data[4,11] = 0.3990511495040125
$data[4,12] = 19.87807525522901
$data[4,13] = 178.902677297061
$data[4,14] = 0.1612247718381573
$data[4,15] = 0.1612247718381573

$data[5,0] = 3
$data[5,1] = 1
$data[5,2] = 11.94315233333334
$data[5,3] = 35.829457
$data[5,4] = 0.4040203167903319
$data[5,5] = 0.4040203167903319
$data[5,6] = 3
$data[5,7] = 1
$data[5,8] = 0.437958
$data[5,9] = 1.313874
$data[5,10] = 0.1050039584054939
$data[5,11] = 0.1050039584054938
$data[5,12] = 5.230599109602001
$data[5,13] = 47.075391986418
$data[5,14] = 0.04242373253922646
$data[5,15] = 0.04242373253922645

$data[6,0] = 3
$data[6,1] = 1
$data[6,2] = 11.94315233333334
$data[6,3] = 35.829457
$data[6,4] = 0.4040203167903319
$data[6,5] = 0.4040203167903319
$data[6,6] = 3
$data[6,7] = 1
$data[6,8] = 1.439215333333333
$data[6,9] = 4.317646
$data[6,10] = 0.3450634695516061
$data[6,11] = 0.3450634695516061
$data[6,12] = 17.18876796646911
$data[6,13] = 154.698911698222
$data[6,14] = 0.1394126522810109
$data[6,15] = 0.1394126522810109

$data[7,0] = 3
$data[7,1] = 1
$data[7,2] = 11.94315233333334
$data[7,3] = 35.829457
$data[7,4] = 0.4040203167903319
$data[7,5] = 0.4040203167903319
$data[7,6] = 3
$data[7,7] = 1
$data[7,8] = 0.629307
$data[7,9] = 1.887921
$data[7,10] = 0.1508814225388875
$data[7,11] = 0.1508814225388875
$data[7,12] = 7.515909365433001
$data[7,13] = 67.643184288897
$data[7,14] = 0.06095916013193727
$data[7,15] = 0.06095916013193727

$data[8,0] = 3
$data[8,1] = 1
$data[8,2] = 12.844759
$data[8,3] = 38.534277
$data[8,4] = 0.4345204227020912
$data[8,5] = 0.4345204227020912
$data[8,6] = 3
$data[8,7] = 1
$data[8,8] = 1.664391
$data[8,9] = 4.993173000000001
$data[8,10] = 0.3990511495040125
$data[8,11] = 0.3990511495040125
$data[8,12] = 21.378701276769
$data[8,13] = 192.408311490921
$data[8,14] = 0.1733958741622389
$data[8,15] = 0.1733958741622389

$data[9,0] = 3
$data[9,1] = 1
$data[9,2] = 12.844759
$data[9,3] = 38.534277
$data[9,4] = 0.4345204227020912
$data[9,5] = 0.4345204227020912
$data[9,6] = 3
$data[9,7] = 1
$data[9,8] = 0.437958
$data[9,9] = 1.313874
$data[9,10] = 0.1050039584054939
$data[9,11] = 0.1050039584054938
$data[9,12] = 5.625464962122001
$data[9,13] = 50.629184659098
$data[9,14] = 0.04562636439174799
$data[9,15] = 0.04562636439174798

$data[10,0] = 3
$data[10,1] = 1
$data[10,2] = 12.844759
$data[10,3] = 38.534277
$data[10,4] = 0.4345204227020912
$data[10,5] = 0.4345204227020912
$data[10,6] = 3
$data[10,7] = 1
$data[10,8] = 1.439215333333333
$data[10,9] = 4.317646
$data[10,10] = 0.3450634695516061
$data[10,11] = 0.3450634695516061
$data[10,12] = 18.48637410577134
$data[10,13] = 166.377366951942
$data[10,14] = 0.149937124648614
$data[10,15] = 0.149937124648614

$data[11,0] = 3
$data[11,1] = 1
$data[11,2] = 12.844759
$data[11,3] = 38.534277
$data[11,4] = 0.4345204227020912
$data[11,5] = 0.4345204227020912
$data[11,6] = 3
$data[11,7] = 1
$data[11,8] = 0.629307
$data[11,9] = 1.887921
$data[11,10] = 0.1508814225388875
$data[11,11] = 0.1508814225388875
$data[11,12] = 8.083296752013
$data[11,13] = 72.749670768117
$data[11,14] = 0.06556105949949023
$data[11,15] = 0.06556105949949023

$data[12,0] = 3
$data[12,1] = 1
$data[12,2] = 4.313537
$data[12,3] = 12.940611
$data[12,4] = 0.1459209877414679
$data[12,5] = 0.145920987741468
$data[12,6] = 3
$data[12,7] = 1
$data[12,8] = 1.664391
$data[12,9] = 4.993173000000001
$data[12,10] = 0.3990511495040125
$data[12,11] = 0.3990511495040125
$data[12,12] = 7.179412160967001
$data[12,13] = 64.61470944870301
$data[12,14] = 0.05822993789499371
$data[12,15] = 0.05822993789499372

$data[13,0] = 3
$data[13,1] = 1
$data[13,2] = 4.313537
$data[13,3] = 12.940611
$data[13,4] = 0.1459209877414679
$data[13,5] = 0.145920987741468
$data[13,6] = 3
$data[13,7] = 1
$data[13,8] = 0.437958
$data[13,9] = 1.313874
$data[13,10] = 0.1050039584054939
$data[13,11] = 0.1050039584054938
$data[13,12] = 1.889148037446
$data[13,13] = 17.002332337014
$data[13,14] = 0.01532228132729368
$data[13,15] = 0.01532228132729368

$data[14,0] = 3
$data[14,1] = 1
$data[14,2] = 4.313537
$data[14,3] = 12.940611
$data[14,4] = 0.1459209877414679
$data[14,5] = 0.145920987741468
$data[14,6] = 3
$data[14,7] = 1
$data[14,8] = 1.439215333333333
$data[14,9] = 4.317646
$data[14,10] = 0.3450634695516061
$data[14,11] = 0.3450634695516061
$data[14,12] = 6.208108591300666
$data[14,13] = 55.872977321706
$data[14,14] = 0.0503520023104683
$data[14,15] = 0.05035200231046831

$data[15,0] = 3
$data[15,1] = 1
$data[15,2] = 4.313537
$data[15,3] = 12.940611
$data[15,4] = 0.1459209877414679
$data[15,5] = 0.145920987741468
$data[15,6] = 3
$data[15,7] = 1
$data[15,8] = 0.629307
$data[15,9] = 1.887921
$data[15,10] = 0.1508814225388875
$data[15,11] = 0.1508814225388875
$data[15,12] = 2.714539028859
$data[15,13] = 24.430851259731
$data[15,14] = 0.02201676620871225
$data[15,15] = 0.02201676620871226

$ws.Range("E2:T17").Value = $data
Write-Output "done"
